$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 29 hold the "Förändrad" date value.
# Update the date serial value from 45174 (2023-09-05) to 45175 (2023-09-06).
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Range("C$row")
    $current = $cell.Value()
    if ($current.ToOADate() -eq 45174) {
        $cell.Value = [DateTime]::FromOADate(45175)
    }
}
